$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row,Fecha(D),Volumen(J),PrecioMinimo(K),PrecioMaximo(L),PrecioPromedioPonderado(M),PrecioPorKg(P)
$rows = @(
    "2,44599,400,15000,16000,15500,1192",
    "3,44309,400,26000,27000,26500,2038",
    "4,44435,480,13000,14000,13500,1038",
    "5,44414,500,14000,15000,14500,1115",
    "6,44657,460,15000,16000,15500,1192",
    "7,44426,460,14000,15000,14500,1115",
    "8,44428,480,14000,15000,14500,1115",
    "9,44680,400,13500,14000,13750,1058",
    "10,44670,480,14500,15000,14750,1135",
    "11,44484,360,14000,15000,14500,1115",
    "12,44312,400,26000,27000,26500,2038",
    "13,44379,600,17000,18000,17500,1346",
    "14,44419,600,14000,15000,14500,1115",
    "15,44596,500,16000,17000,16500,1269",
    "16,44687,440,14000,15000,14500,1115",
    "17,44260,400,37000,38000,37500,2885",
    "18,44383,200,17000,18000,17500,1346",
    "19,44644,400,15000,16000,15500,1192",
    "20,44335,480,24500,25000,24750,1904",
    "21,44418,500,14000,15000,14500,1115",
    "22,44410,600,14000,15000,14500,1115",
    "23,44582,520,15000,16000,15500,1192",
    "24,44333,440,24000,25000,24500,1885",
    "25,44498,400,14000,15000,14500,1115",
    "26,44658,400,15000,16000,15500,1192",
    "27,44400,600,15000,16000,15500,1192",
    "28,44412,600,14000,15000,14500,1115",
    "29,44694,400,13000,14000,13500,1038",
    "30,44505,400,16000,17000,16500,1269",
    "31,44631,400,16000,17000,16500,1269",
    "32,44365,500,19500,20000,19750,1519",
    "33,44575,500,14000,15000,14500,1115",
    "34,44445,600,13000,14000,13500,1038",
    "35,44533,520,17000,18000,17500,1346",
    "36,44344,400,18500,19000,18750,1442",
    "37,44323,460,25000,26000,25500,1962",
    "38,44326,460,25000,26000,25500,1962",
    "39,44692,400,14000,15000,14500,1115",
    "40,44442,460,14000,15000,14500,1115"
)

foreach ($line in $rows) {
    $parts = $line.Split(",")
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 4).Value = [double]$parts[1]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = [double]$parts[2]  # J - Volumen
    $ws.Cells.Item($r, 11).Value = [double]$parts[3]  # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = [double]$parts[4]  # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = [double]$parts[5]  # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = [double]$parts[6]  # P - Precio $/Kg
}
